$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 320, shifting the existing rows 320:383 down to 321:384.
$ws.Rows.Item(320).Insert()

# Populate the newly inserted row 320 with the new weekly observation.
$ws.Range("A320").Value = 7
$ws.Range("B320").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C320").Value = "Ñuble"
$ws.Range("D320").Value = 45211
$ws.Range("E320").Value = 16
$ws.Range("F320").Value = 100112043
$ws.Range("G320").Value = "Pepino ensalada"
$ws.Range("H320").Value = "Sin especificar"
$ws.Range("I320").Value = "Primera"
$ws.Range("J320").Value = 100
$ws.Range("K320").Value = 14500
$ws.Range("L320").Value = 15000
$ws.Range("M320").Value = 14800
$ws.Range("N320").Value = "`$/caja 60 unidades"
$ws.Range("O320").Value = "Región de Arica y Parinacota"
$ws.Range("P320").Value = 247
$ws.Range("Q320").Value = 60
$ws.Range("R320").Value = "Hortaliza"
